# Refactor the export template: the header row gains two new leading
# columns ("Name" and "Age") while the pre-existing columns become
# "Email" (now column C) and "Games" (now column D). The "Games" list
# data-validation dropdown that used to sit on column B follows the
# "Games" header and now applies to column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header layout for row 1.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Games"

# Move the "Games" list validation from column B to column D.
$ws.Range("B2:B100000").Validation.Delete()

$validation = $ws.Range("D2:D100000").Validation
$validation.Add(3, 1, 1, '"Super Mario,SONIC,Zelda,GTA"')
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $false
$validation.ShowError = $false
